$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price / Volume(1h) columns are stored as plain text (values like
# "261.59" or "0.89%"), not real numbers. Mark each touched cell as Text
# first so Excel does not silently reinterpret the refreshed values as
# numeric/percentage, then write the new scraped value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.89%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.47%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.709"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.27%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06207"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.90%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.729"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.86%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8502"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.12%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9094"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.43%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1409"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.78%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04678"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.72%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.23%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03152"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.07%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09060"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.73%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001536"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.51%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006167"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.58%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006128"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.78%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.470"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.11%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.170"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.04%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.62%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.88%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.085"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.86%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04232"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.28%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.67%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.36%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.13%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03901"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.41%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.07%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004129"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002162"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.72%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01349"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.65%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005173"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.59%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.05%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.03590"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-34.15%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1667"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "26.16%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.05%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
